# Auto-generated Excel COM-interop script
# Updates 'want to go' counts (column F) and minimum price (column G)
# on the four sheets to match gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1928
$ws.Range("F5").Value = 91
$ws.Range("F7").Value = 1795
$ws.Range("F9").Value = 1247
$ws.Range("F10").Value = 64
$ws.Range("F11").Value = 364
$ws.Range("F13").Value = 2516
$ws.Range("F14").Value = 344
$ws.Range("F15").Value = 856
$ws.Range("F16").Value = 1060
$ws.Range("F18").Value = 49
$ws.Range("F19").Value = 1512
$ws.Range("F20").Value = 408485
$ws.Range("F22").Value = 152
$ws.Range("F25").Value = 1358
$ws.Range("F26").Value = 942
$ws.Range("F27").Value = 45
$ws.Range("F29").Value = 179
$ws.Range("F31").Value = 413
$ws.Range("F35").Value = 1803
$ws.Range("F36").Value = 436
$ws.Range("F37").Value = 30
$ws.Range("F38").Value = 144
$ws.Range("F40").Value = 2213
$ws.Range("F42").Value = 873
$ws.Range("F43").Value = 2709
$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 28
$ws.Range("F12").Value = 354
$ws.Range("F14").Value = 42
$ws.Range("F19").Value = 209
$ws.Range("F22").Value = 265
$ws.Range("F26").Value = 51
$ws.Range("G26").Value = 190
$ws.Range("F27").Value = 51
$ws.Range("G27").Value = 190
$ws.Range("F36").Value = 157
$ws = $wb.Worksheets.Item(3)
$ws.Range("F6").Value = 4822
$ws.Range("F10").Value = 901
$ws.Range("F12").Value = 568
$ws.Range("F13").Value = 1304
$ws.Range("F14").Value = 369
$ws.Range("F15").Value = 1092
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1928
$ws.Range("F5").Value = 4822
$ws.Range("F7").Value = 901
$ws.Range("F9").Value = 568
$ws.Range("F10").Value = 1304
$ws.Range("F12").Value = 1795
$ws.Range("F14").Value = 1247
$ws.Range("F16").Value = 364
$ws.Range("F17").Value = 1092
$ws.Range("F18").Value = 2516
$ws.Range("F19").Value = 28
$ws.Range("F20").Value = 344
$ws.Range("F21").Value = 856
$ws.Range("F22").Value = 1060
$ws.Range("F23").Value = 571
$ws.Range("F24").Value = 1512
$ws.Range("F26").Value = 354
$ws.Range("F28").Value = 152
$ws.Range("F29").Value = 588
$ws.Range("F30").Value = 1358
$ws.Range("F31").Value = 942
$ws.Range("F33").Value = 179
$ws.Range("F37").Value = 413
$ws.Range("F40").Value = 1803
$ws.Range("F41").Value = 51
$ws.Range("G41").Value = 190
$ws.Range("F42").Value = 30
$ws.Range("F43").Value = 144
$ws.Range("F44").Value = 2213
$ws.Range("F46").Value = 873
$ws.Range("F47").Value = 2709
